$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder "Periodo Mora" column (E16:E28) so periods run chronologically
# ascending (2012, 2101, 2102, ... 2112) instead of the old descending order.
$periods = @("2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# The "Valor Mora" outlier (21333 instead of 40000) now belongs with period
# 2112, which after the reorder sits in row 28 instead of row 16.
$ws.Range("F16").Value = 40000
$ws.Range("F28").Value = 21333

$wb.Save()
